$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (right after the header), shifting existing
# data rows down by one.
$ws.Rows.Item(2).Insert()

# Copy the formatting used by the rest of column A (border + center
# alignment) down into the newly-inserted cell, then set its values.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A2").Value = "CA"
$ws.Range("B2").Value = 1996

# The data range used to sort the table moved down by one row along with
# the data (it still excludes the newly-inserted row).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A3:A43"))
$ws.Sort.SetRange($ws.Range("A3:B43"))
$ws.Sort.Header = -4142  # xlNo
$ws.Sort.Apply()

# Update the selection to match the post-edit state.
$ws.Range("B3").Select()
